$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 79715
$ws.Range("B2").Value = "Oliver Cavalcanti"
$ws.Range("D2").Value = "Doenca"
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 45088
$ws.Range("G2").Value = 8902.870000000001

# Row 3
$ws.Range("A3").Value = 54487
$ws.Range("B3").Value = "Miguel Camargo"
$ws.Range("C3").Value = "P&D"
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 45094
$ws.Range("G3").Value = 5996.19

# Row 4
$ws.Range("A4").Value = 34073
$ws.Range("B4").Value = "Lunna Ramos"
$ws.Range("C4").Value = "P&D"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 45105
$ws.Range("G4").Value = 9808.66

# Row 5
$ws.Range("A5").Value = 18482
$ws.Range("B5").Value = "Ravy Costela"
$ws.Range("C5").Value = "Atendimento ao Cliente"
$ws.Range("D5").Value = "Doenca"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 45102
$ws.Range("G5").Value = 7328.48

# Row 6
$ws.Range("A6").Value = 98317
$ws.Range("B6").Value = "Ana Júlia Brito"
$ws.Range("C6").Value = "TI"
$ws.Range("D6").Value = "Consulta medica"
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 45089
$ws.Range("G6").Value = 7284.62

# Row 7
$ws.Range("A7").Value = 11661
$ws.Range("B7").Value = "Eloah Duarte"
$ws.Range("C7").Value = "Marketing"
$ws.Range("D7").Value = "Outros"
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 45095
$ws.Range("G7").Value = 7987.61

# Row 8
$ws.Range("A8").Value = 17971
$ws.Range("B8").Value = "Letícia Martins"
$ws.Range("C8").Value = "Vendas"
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 45088
$ws.Range("G8").Value = 4678.77

# Row 9
$ws.Range("A9").Value = 95524
$ws.Range("B9").Value = "Srta. Luiza Costa"
$ws.Range("C9").Value = "Vendas"
$ws.Range("D9").Value = "Viagem de negocios"
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = 45102
$ws.Range("G9").Value = 4403.61

# Row 10
$ws.Range("A10").Value = 31766
$ws.Range("B10").Value = "Ryan Pereira"
$ws.Range("C10").Value = "TI"
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 45094
$ws.Range("G10").Value = 5454.03

# Row 11
$ws.Range("A11").Value = 54385
$ws.Range("B11").Value = "Valentina Cirino"
$ws.Range("C11").Value = "P&D"
$ws.Range("D11").Value = "Viagem de negocios"
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 45089
$ws.Range("G11").Value = 5594.65
